# Generate Report for Handback
# Rename the handed-back GUID-named files and refresh the correspond
# handoff/handback xlf file names + timestamps.

$wb = $excel.ActiveWorkbook

$md1Old = "a8a25305-d45e-4370-9a91-f732f20bbebb.md"
$md2Old = "cd7831bf-1df1-47ac-b351-225d4c84c32e.md"

$md1New = "8d8d2085-fe47-481b-b662-0dc702c42ed3.md"
$md2New = "ffff398fff67-c961-4876-b775-6cd7871573cb.md"

$xlfZhNew = "8d8d2085-fe47-481b-b662-0dc702c42ed3.bc702ab10f0cf6487cf672ab2645395a00ab4626.zh-cn.xlf"
$xlfDeNew = "8d8d2085-fe47-481b-b662-0dc702c42ed3.bc702ab10f0cf6487cf672ab2645395a00ab4626.de-de.xlf"

$tZhHandoff = "2016-03-22 15:11:02"
$tZhHandback = "2016-03-22 15:11:26"
$tDeHandoff = "2016-03-22 15:11:09"
$tDeHandback = "2016-03-22 15:11:36"

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$addr,
        [string]$newValue
    )
    $ws.Range($addr).Value = $newValue
    $target = $ws.Range($addr).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $newValue
        }
    }
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $wsOverview "A2" $md1New
Set-CellAndHyperlink $wsOverview "A3" $md2New

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $wsZh "A2" $md1New
Set-CellAndHyperlink $wsZh "D2" $xlfZhNew
Set-CellAndHyperlink $wsZh "E2" $tZhHandoff
Set-CellAndHyperlink $wsZh "F2" $md1New
Set-CellAndHyperlink $wsZh "G2" $xlfZhNew
Set-CellAndHyperlink $wsZh "H2" $tZhHandback

Set-CellAndHyperlink $wsZh "A3" $md2New
Set-CellAndHyperlink $wsZh "D3" $xlfZhNew
Set-CellAndHyperlink $wsZh "E3" $tZhHandoff
Set-CellAndHyperlink $wsZh "F3" $md2New
Set-CellAndHyperlink $wsZh "G3" $xlfZhNew
Set-CellAndHyperlink $wsZh "H3" $tZhHandback

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $wsDe "A2" $md1New
Set-CellAndHyperlink $wsDe "D2" $xlfDeNew
Set-CellAndHyperlink $wsDe "E2" $tDeHandoff
Set-CellAndHyperlink $wsDe "F2" $md1New
Set-CellAndHyperlink $wsDe "G2" $xlfDeNew
Set-CellAndHyperlink $wsDe "H2" $tDeHandback

Set-CellAndHyperlink $wsDe "A3" $md2New
Set-CellAndHyperlink $wsDe "D3" $xlfDeNew
Set-CellAndHyperlink $wsDe "E3" $tDeHandoff
Set-CellAndHyperlink $wsDe "F3" $md2New
Set-CellAndHyperlink $wsDe "G3" $xlfDeNew
Set-CellAndHyperlink $wsDe "H3" $tDeHandback

Write-Output "Handback report regenerated."
